$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 146.17053
$ws.Range("H2").Value = 438.51159
$ws.Range("I2").Value = 0.4047435297111188
$ws.Range("J2").Value = 0.4047435297111188
$ws.Range("O2").Value = 0.01611173663836548
$ws.Range("P2").Value = 0.01611173663836548
$ws.Range("Q2").Value = 9.311208931529999
$ws.Range("R2").Value = 83.80088038376999
$ws.Range("S2").Value = 0.006521121156788
$ws.Range("T2").Value = 0.006521121156787998
$ws.Range("G3").Value = 146.17053
$ws.Range("H3").Value = 438.51159
$ws.Range("I3").Value = 0.4047435297111188
$ws.Range("J3").Value = 0.4047435297111188
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.701496333333333
$ws.Range("N3").Value = 8.104489000000001
$ws.Range("O3").Value = 0.68328279700753
$ws.Range("P3").Value = 0.68328279700753
$ws.Range("Q3").Value = 394.87915083639
$ws.Range("R3").Value = 3553.91235752751
$ws.Range("S3").Value = 0.2765542910517136
$ws.Range("T3").Value = 0.2765542910517136
$ws.Range("G4").Value = 146.17053
$ws.Range("H4").Value = 438.51159
$ws.Range("I4").Value = 0.4047435297111188
$ws.Range("J4").Value = 0.4047435297111188
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.188504333333333
$ws.Range("N4").Value = 3.565513
$ws.Range("O4").Value = 0.3006054663541045
$ws.Range("P4").Value = 0.3006054663541044
$ws.Range("Q4").Value = 173.72430831063
$ws.Range("R4").Value = 1563.51877479567
$ws.Range("S4").Value = 0.1216681175026172
$ws.Range("T4").Value = 0.1216681175026172
$ws.Range("H5").Value = 632.3552549999999
$ws.Range("I5").Value = 0.5836600531814327
$ws.Range("J5").Value = 0.5836600531814327
$ws.Range("O5").Value = 0.01611173663836548
$ws.Range("P5").Value = 0.01611173663836548
$ws.Range("S5").Value = 0.009403777063193633
$ws.Range("T5").Value = 0.009403777063193631
$ws.Range("H6").Value = 632.3552549999999
$ws.Range("I6").Value = 0.5836600531814327
$ws.Range("J6").Value = 0.5836600531814327
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.701496333333333
$ws.Range("N6").Value = 8.104489000000001
$ws.Range("O6").Value = 0.68328279700753
$ws.Range("P6").Value = 0.68328279700753
$ws.Range("Q6").Value = 569.435134248855
$ws.Range("R6").Value = 5124.916208239695
$ws.Range("S6").Value = 0.3988048736393731
$ws.Range("T6").Value = 0.3988048736393731
$ws.Range("H7").Value = 632.3552549999999
$ws.Range("I7").Value = 0.5836600531814327
$ws.Range("J7").Value = 0.5836600531814327
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.188504333333333
$ws.Range("N7").Value = 3.565513
$ws.Range("O7").Value = 0.3006054663541045
$ws.Range("P7").Value = 0.3006054663541044
$ws.Range("Q7").Value = 250.518986924535
$ws.Range("R7").Value = 2254.670882320815
$ws.Range("S7").Value = 0.175451402478866
$ws.Range("T7").Value = 0.175451402478866
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.5494936666666667
$ws.Range("H8").Value = 1.648481
$ws.Range("I8").Value = 0.00152153793381314
$ws.Range("J8").Value = 0.00152153793381314
$ws.Range("O8").Value = 0.01611173663836548
$ws.Range("P8").Value = 0.01611173663836548
$ws.Range("Q8").Value = 0.03500329606033333
$ws.Range("R8").Value = 0.315029664543
$ws.Range("S8").Value = [double]"2.451461847488008E-05"
$ws.Range("T8").Value = [double]"2.451461847488008E-05"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.5494936666666667
$ws.Range("H9").Value = 1.648481
$ws.Range("I9").Value = 0.00152153793381314
$ws.Range("J9").Value = 0.00152153793381314
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.701496333333333
$ws.Range("N9").Value = 8.104489000000001
$ws.Range("O9").Value = 0.68328279700753
$ws.Range("P9").Value = 0.68328279700753
$ws.Range("Q9").Value = 1.484455125689889
$ws.Range("R9").Value = 13.360096131209
$ws.Range("S9").Value = 0.001039640695168901
$ws.Range("T9").Value = 0.001039640695168901
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5494936666666667
$ws.Range("H10").Value = 1.648481
$ws.Range("I10").Value = 0.00152153793381314
$ws.Range("J10").Value = 0.00152153793381314
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.188504333333333
$ws.Range("N10").Value = 3.565513
$ws.Range("O10").Value = 0.3006054663541045
$ws.Range("P10").Value = 0.3006054663541044
$ws.Range("Q10").Value = 0.6530756039725556
$ws.Range("R10").Value = 5.877680435753001
$ws.Range("S10").Value = 0.0004573826201693595
$ws.Range("T10").Value = 0.0004573826201693594
$ws.Range("G11").Value = 3.410044
$ws.Range("H11").Value = 10.230132
$ws.Range("I11").Value = 0.009442349596941478
$ws.Range("J11").Value = 0.009442349596941478
$ws.Range("O11").Value = 0.01611173663836548
$ws.Range("P11").Value = 0.01611173663836548
$ws.Range("Q11").Value = 0.217223212844
$ws.Range("R11").Value = 1.955008915596
$ws.Range("S11").Value = 0.0001521326499532975
$ws.Range("T11").Value = 0.0001521326499532975
$ws.Range("G12").Value = 3.410044
$ws.Range("H12").Value = 10.230132
$ws.Range("I12").Value = 0.009442349596941478
$ws.Range("J12").Value = 0.009442349596941478
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.701496333333333
$ws.Range("N12").Value = 8.104489000000001
$ws.Range("O12").Value = 0.68328279700753
$ws.Range("P12").Value = 0.68328279700753
$ws.Range("Q12").Value = 9.212221362505334
$ws.Range("R12").Value = 82.909992262548
$ws.Range("S12").Value = 0.006451795042921097
$ws.Range("T12").Value = 0.006451795042921097
$ws.Range("G13").Value = 3.410044
$ws.Range("H13").Value = 10.230132
$ws.Range("I13").Value = 0.009442349596941478
$ws.Range("J13").Value = 0.009442349596941478
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.188504333333333
$ws.Range("N13").Value = 3.565513
$ws.Range("O13").Value = 0.3006054663541045
$ws.Range("P13").Value = 0.3006054663541044
$ws.Range("Q13").Value = 4.052852070857333
$ws.Range("R13").Value = 36.475668637716
$ws.Range("S13").Value = 0.002838421904067083
$ws.Range("T13").Value = 0.002838421904067083
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.228434
$ws.Range("H14").Value = 0.6853020000000001
$ws.Range("I14").Value = 0.0006325295766939459
$ws.Range("J14").Value = 0.0006325295766939459
$ws.Range("O14").Value = 0.01611173663836548
$ws.Range("P14").Value = 0.01611173663836548
$ws.Range("Q14").Value = 0.014551474234
$ws.Range("R14").Value = 0.130963268106
$ws.Range("S14").Value = [double]"1.019114995566966E-05"
$ws.Range("T14").Value = [double]"1.019114995566965E-05"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.228434
$ws.Range("H15").Value = 0.6853020000000001
$ws.Range("I15").Value = 0.0006325295766939459
$ws.Range("J15").Value = 0.0006325295766939459
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 2.701496333333333
$ws.Range("N15").Value = 8.104489000000001
$ws.Range("O15").Value = 0.68328279700753
$ws.Range("P15").Value = 0.68328279700753
$ws.Range("Q15").Value = 0.6171136134086668
$ws.Range("R15").Value = 5.554022520678001
$ws.Range("S15").Value = 0.0004321965783534283
$ws.Range("T15").Value = 0.0004321965783534283
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.228434
$ws.Range("H16").Value = 0.6853020000000001
$ws.Range("I16").Value = 0.0006325295766939459
$ws.Range("J16").Value = 0.0006325295766939459
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.188504333333333
$ws.Range("N16").Value = 3.565513
$ws.Range("O16").Value = 0.3006054663541045
$ws.Range("P16").Value = 0.3006054663541044
$ws.Range("Q16").Value = 0.2714947988806667
$ws.Range("R16").Value = 2.443453189926001
$ws.Range("S16").Value = 0.0001901418483848479
$ws.Range("T16").Value = 0.0001901418483848479
